$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row-level "puesto" (column C) corrections
$ws.Range("C2").Value  = "Metrologo Jr."
$ws.Range("C5").Value  = "Metrologo"
$ws.Range("C7").Value  = "Metrologo Jr."
$ws.Range("C8").Value  = "Metrologo"
$ws.Range("C9").Value  = "Metrologo Jr."
$ws.Range("C11").Value = "Metrologo"
$ws.Range("C13").Value = "Metrologo Jr."
$ws.Range("C14").Value = "Metrologo"
$ws.Range("C18").Value = "Metrologo Jr."
$ws.Range("C24").Value = "Metrologo Jr."
$ws.Range("C25").Value = "Metrologo Jr."

# Salary correction
$ws.Range("I16").Value = 108000

# Restore selection to a single cell
$ws.Range("E6").Select()
